$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M70").Value = -2311.33338
$ws.Range("N70").Value = -4906.6155
$ws.Range("I70").Value = 860.44446
$ws.Range("L70").Value = 4366.6155
$ws.Range("K70").Value = 2581.33338
$ws.Range("H70").Value = 1212.091
$ws.Range("J70").Value = 1455.5385
$ws.Range("N73").Value = -6238.6155
$ws.Range("I73").Value = 860.44446
$ws.Range("J73").Value = 1455.5385
$ws.Range("M73").Value = -1645.33338
$ws.Range("K73").Value = 2581.33338
$ws.Range("L73").Value = 4366.6155
$ws.Range("H73").Value = 1212.091
$ws.Range("K80").Value = 1039.59999
$ws.Range("H80").Value = 1273.579
$ws.Range("L80").Value = 14250
$ws.Range("I80").Value = 346.53333
$ws.Range("J80").Value = 4750
$ws.Range("N80").Value = -16246
$ws.Range("M80").Value = -41.59998999999993
$ws.Range("K83").Value = 3118.79997
$ws.Range("M83").Value = 1873.20003
$ws.Range("N83").Value = -52734
$ws.Range("J83").Value = 4750
$ws.Range("L83").Value = 42750
$ws.Range("I83").Value = 346.53333
$ws.Range("H83").Value = 1273.579
$ws.Range("H132").Value = 36589.207
$ws.Range("K132").Value = 113473.605
$ws.Range("M132").Value = -110943.605
$ws.Range("I132").Value = 37824.535
$ws.Range("N133").Value = -56236
$ws.Range("L133").Value = 46116
$ws.Range("H133").Value = 46116
$ws.Range("J133").Value = 46116

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("K63").Value = 10671.6
$ws.Range("I63").Value = 10671.6
$ws.Range("H63").Value = 9208.286
$ws.Range("M63").Value = -9985.6
$ws.Range("I66").Value = 10671.6
$ws.Range("K66").Value = 53358
$ws.Range("M66").Value = -49926
$ws.Range("H66").Value = 9208.286
$ws.Range("J75").Value = 49000
$ws.Range("N75").Value = -50748
$ws.Range("L75").Value = 49000
$ws.Range("H75").Value = 49000
$ws.Range("L78").Value = 147000
$ws.Range("J78").Value = 49000
$ws.Range("H78").Value = 49000
$ws.Range("N78").Value = -155736
$ws.Range("N133").Value = -50625.25
$ws.Range("L133").Value = 45565.25
$ws.Range("H133").Value = 45565.25
$ws.Range("J133").Value = 45565.25
$ws.Range("L139").Value = 51715
$ws.Range("J139").Value = 51715
$ws.Range("H139").Value = 51715
$ws.Range("N139").Value = -61995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("J59").Value = 54000
$ws.Range("N59").Value = -55694
$ws.Range("L59").Value = 54000
$ws.Range("H59").Value = 54000
$ws.Range("N133").Value = -50900
$ws.Range("L133").Value = 40780
$ws.Range("H133").Value = 40780
$ws.Range("J133").Value = 40780
$ws.Range("M134").Value = -8785.9095
$ws.Range("I134").Value = 3773.6365
$ws.Range("N134").Value = -28154.181
$ws.Range("H134").Value = 5734.1816
$ws.Range("L134").Value = 23084.181
$ws.Range("J134").Value = 7694.727
$ws.Range("K134").Value = 11320.9095

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 1212.25
$ws.Range("N16").Value = -1786.25
$ws.Range("I16").Value = 889
$ws.Range("H16").Value = 988.46155
$ws.Range("M16").Value = -602
$ws.Range("K16").Value = 889
$ws.Range("L16").Value = 1212.25
$ws.Range("I31").Value = 1871.963
$ws.Range("K31").Value = 1871.963
$ws.Range("H31").Value = 6258.421
$ws.Range("M31").Value = -1576.963
$ws.Range("H34").Value = 6258.421
$ws.Range("K34").Value = 1871.963
$ws.Range("I34").Value = 1871.963
$ws.Range("M34").Value = -1669.963
$ws.Range("N51").Value = -21471
$ws.Range("H51").Value = 19999
$ws.Range("L51").Value = 19999
$ws.Range("J51").Value = 19999
$ws.Range("N61").Value = -20695
$ws.Range("L61").Value = 19999
$ws.Range("J61").Value = 19999
$ws.Range("H61").Value = 19999
$ws.Range("H80").Value = 93213
$ws.Range("L80").Value = 93213
$ws.Range("J80").Value = 93213
$ws.Range("N80").Value = -95459
$ws.Range("H81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("N83").Value = -290871
$ws.Range("J83").Value = 93213
$ws.Range("L83").Value = 279639
$ws.Range("H83").Value = 93213
$ws.Range("J84").Value = 0
$ws.Range("H84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("H86").Value = 50002000
$ws.Range("M86").Value = -71428587
$ws.Range("N86").Value = -6245.6667
$ws.Range("I86").Value = 71429710
$ws.Range("K86").Value = 71429710
$ws.Range("J86").Value = 3999.6667
$ws.Range("L86").Value = 3999.6667
$ws.Range("N89").Value = -31230.3335
$ws.Range("H89").Value = 50002000
$ws.Range("L89").Value = 19998.3335
$ws.Range("K89").Value = 357148550
$ws.Range("J89").Value = 3999.6667
$ws.Range("M89").Value = -357142934
$ws.Range("I89").Value = 71429710
$ws.Range("H94").Value = 2098.8
$ws.Range("L94").Value = 2272.5
$ws.Range("J94").Value = 2272.5
$ws.Range("I94").Value = 1404
$ws.Range("N94").Value = -3174.5
$ws.Range("K94").Value = 1404
$ws.Range("M94").Value = -953
$ws.Range("L99").Value = 3637.375
$ws.Range("N99").Value = -6633.375
$ws.Range("I99").Value = 1816.1428
$ws.Range("M99").Value = -318.1428000000001
$ws.Range("H99").Value = 2787.4666
$ws.Range("J99").Value = 3637.375
$ws.Range("K99").Value = 1816.1428
$ws.Range("M107").Value = 1628.27274
$ws.Range("N107").Value = -4357
$ws.Range("I107").Value = 291.72726
$ws.Range("K107").Value = 291.72726
$ws.Range("L107").Value = 517
$ws.Range("J107").Value = 517
$ws.Range("H107").Value = 326.3846
$ws.Range("I113").Value = 889
$ws.Range("L113").Value = 1212.25
$ws.Range("N113").Value = -5552.25
$ws.Range("H113").Value = 988.46155
$ws.Range("M113").Value = 1281
$ws.Range("J113").Value = 1212.25
$ws.Range("K113").Value = 889
$ws.Range("J122").Value = 1714.2858
$ws.Range("H122").Value = 1781
$ws.Range("I122").Value = 1936.6666
$ws.Range("L122").Value = 5142.857400000001
$ws.Range("M122").Value = -3359.9998
$ws.Range("N122").Value = -10042.8574
$ws.Range("K122").Value = 5809.9998
$ws.Range("J126").Value = 3637.375
$ws.Range("M126").Value = -2978.428400000001
$ws.Range("L126").Value = 10912.125
$ws.Range("I126").Value = 1816.1428
$ws.Range("H126").Value = 2787.4666
$ws.Range("N126").Value = -15852.125
$ws.Range("K126").Value = 5448.428400000001
$ws.Range("H132").Value = 2924.5715
$ws.Range("K132").Value = 5052
$ws.Range("M132").Value = -2522
$ws.Range("N132").Value = -16625
$ws.Range("I132").Value = 1684
$ws.Range("L132").Value = 11565
$ws.Range("J132").Value = 3855
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("K61").Value = 510.66669
$ws.Range("I61").Value = 170.22223
$ws.Range("H61").Value = 170.22223
$ws.Range("M61").Value = -295.66669
$ws.Range("J122").Value = 1039.8182
$ws.Range("H122").Value = 719.35
$ws.Range("I122").Value = 327.66666
$ws.Range("L122").Value = 9358.363799999999
$ws.Range("M122").Value = -498.9999399999997
$ws.Range("N122").Value = -14258.3638
$ws.Range("K122").Value = 2948.99994

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("N43").Value = -20302
$ws.Range("H43").Value = 6050.6665
$ws.Range("L43").Value = 20000
$ws.Range("I43").Value = 3260.8
$ws.Range("K43").Value = 3260.8
$ws.Range("J43").Value = 20000
$ws.Range("M43").Value = -3109.8
$ws.Range("H137").Value = 46893.332
$ws.Range("L137").Value = 46893.332
$ws.Range("J137").Value = 46893.332
$ws.Range("N137").Value = -57093.332
$ws.Range("L138").Value = 76061.28999999999
$ws.Range("H138").Value = 76061.28999999999
$ws.Range("N138").Value = -86341.28999999999
$ws.Range("J138").Value = 76061.28999999999
$ws.Range("L139").Value = 42061.2
$ws.Range("J139").Value = 42061.2
$ws.Range("H139").Value = 42061.2
$ws.Range("N139").Value = -52341.2

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("K62").Value = 15600.286
$ws.Range("N62").Value = -17459.111
$ws.Range("L62").Value = 16211.111
$ws.Range("M62").Value = -14976.286
$ws.Range("H62").Value = 15839.305
$ws.Range("I62").Value = 15600.286
$ws.Range("J62").Value = 16211.111
$ws.Range("I65").Value = 15600.286
$ws.Range("L65").Value = 81055.55500000001
$ws.Range("H65").Value = 15839.305
$ws.Range("J65").Value = 16211.111
$ws.Range("K65").Value = 78001.42999999999
$ws.Range("N65").Value = -87295.55500000001
$ws.Range("M65").Value = -74881.42999999999
$ws.Range("H132").Value = 4286.968
$ws.Range("K132").Value = 13801.173
$ws.Range("M132").Value = -11271.173
$ws.Range("N132").Value = -15217.625
$ws.Range("I132").Value = 4600.391
$ws.Range("L132").Value = 10157.625
$ws.Range("J132").Value = 3385.875
